$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1955128205128205
$ws.Range("C2").Value = 0.5480769230769231
$ws.Range("J2").Value = 0.02243589743589744
$ws.Range("P2").Value = 0.1474358974358974
$ws.Range("S2").Value = 0.08653846153846154
$ws.Range("B3").Value = 0.02717391304347826
$ws.Range("C3").Value = 0.05978260869565218
$ws.Range("J3").Value = 0.02173913043478261
$ws.Range("O3").Value = 0.005434782608695652
$ws.Range("P3").Value = 0.7065217391304348
$ws.Range("S3").Value = 0.1793478260869565
$ws.Range("J4").Value = 0.05
$ws.Range("P4").Value = 0.575
$ws.Range("S4").Value = 0.375
$ws.Range("P5").Value = 0.9
$ws.Range("S5").Value = 0.1
$ws.Range("B6").Value = 0.0564516129032258
$ws.Range("F6").Value = 0.04032258064516129
$ws.Range("J6").Value = 0.2580645161290323
$ws.Range("O6").Value = 0.0282258064516129
$ws.Range("Q6").Value = 0.1330645161290323
$ws.Range("R6").Value = 0.07258064516129033
$ws.Range("S6").Value = 0.4112903225806452
$ws.Range("B7").Value = 0.08597285067873303
$ws.Range("D7").Value = 0.01357466063348416
$ws.Range("E7").Value = 0.009049773755656109
$ws.Range("F7").Value = 0.04524886877828054
$ws.Range("J7").Value = 0.1538461538461539
$ws.Range("O7").Value = 0.03167420814479638
$ws.Range("Q7").Value = 0.1583710407239819
$ws.Range("R7").Value = 0.1040723981900453
$ws.Range("S7").Value = 0.3981900452488688
$ws.Range("B8").Value = 0.07478632478632478
$ws.Range("D8").Value = 0.01495726495726496
$ws.Range("F8").Value = 0.05982905982905983
$ws.Range("J8").Value = 0.1282051282051282
$ws.Range("O8").Value = 0.02564102564102564
$ws.Range("Q8").Value = 0.1837606837606838
$ws.Range("R8").Value = 0.08974358974358974
$ws.Range("S8").Value = 0.4230769230769231
$ws.Range("B9").Value = 0.1026615969581749
$ws.Range("D9").Value = 0.007604562737642586
$ws.Range("F9").Value = 0.05703422053231939
$ws.Range("J9").Value = 0.1292775665399239
$ws.Range("O9").Value = 0.01520912547528517
$ws.Range("Q9").Value = 0.1711026615969582
$ws.Range("R9").Value = 0.1140684410646388
$ws.Range("S9").Value = 0.403041825095057
$ws.Range("B10").Value = 0.1026878015161957
$ws.Range("D10").Value = 0.01998621640248105
$ws.Range("E10").Value = 0.006202618883528601
$ws.Range("F10").Value = 0.07649896623018608
$ws.Range("J10").Value = 0.1150930392832529
$ws.Range("O10").Value = 0.01929703652653342
$ws.Range("Q10").Value = 0.217091660923501
$ws.Range("R10").Value = 0.0771881461061337
$ws.Range("S10").Value = 0.3659545141281875
$ws.Range("G11").Value = 0.1223880597014925
$ws.Range("J11").Value = 0.08656716417910448
$ws.Range("K11").Value = 0.1582089552238806
$ws.Range("L11").Value = 0.6119402985074627
$ws.Range("S11").Value = 0.0208955223880597
$ws.Range("G12").Value = 0.7201834862385321
$ws.Range("J12").Value = 0.1880733944954129
$ws.Range("K12").Value = 0.004587155963302753
$ws.Range("L12").Value = 0.04128440366972477
$ws.Range("S12").Value = 0.04587155963302753
$ws.Range("G13").Value = 0.725
$ws.Range("J13").Value = 0.175
$ws.Range("S13").Value = 0.1
$ws.Range("G14").Value = 0.6
$ws.Range("J14").Value = 0.4
$ws.Range("F15").Value = 0.01937984496124031
$ws.Range("H15").Value = 0.1550387596899225
$ws.Range("I15").Value = 0.05426356589147287
$ws.Range("J15").Value = 0.3720930232558139
$ws.Range("K15").Value = 0.06589147286821706
$ws.Range("O15").Value = 0.07364341085271318
$ws.Range("S15").Value = 0.2596899224806202
$ws.Range("F16").Value = 0.01970443349753695
$ws.Range("H16").Value = 0.1231527093596059
$ws.Range("I16").Value = 0.1330049261083744
$ws.Range("J16").Value = 0.3497536945812808
$ws.Range("K16").Value = 0.09359605911330049
$ws.Range("M16").Value = 0.03448275862068965
$ws.Range("O16").Value = 0.04926108374384237
$ws.Range("S16").Value = 0.1970443349753695
$ws.Range("F17").Value = 0.0210727969348659
$ws.Range("H17").Value = 0.181992337164751
$ws.Range("I17").Value = 0.1053639846743295
$ws.Range("J17").Value = 0.4061302681992337
$ws.Range("K17").Value = 0.08812260536398467
$ws.Range("M17").Value = 0.01915708812260536
$ws.Range("N17").Value = 0.003831417624521073
$ws.Range("O17").Value = 0.04597701149425287
$ws.Range("S17").Value = 0.1283524904214559
$ws.Range("F18").Value = 0.008733624454148471
$ws.Range("H18").Value = 0.1790393013100437
$ws.Range("I18").Value = 0.1397379912663755
$ws.Range("J18").Value = 0.3755458515283843
$ws.Range("K18").Value = 0.07860262008733625
$ws.Range("M18").Value = 0.01310043668122271
$ws.Range("O18").Value = 0.05676855895196507
$ws.Range("S18").Value = 0.148471615720524
$ws.Range("F19").Value = 0.01499348109517601
$ws.Range("H19").Value = 0.1799217731421121
$ws.Range("I19").Value = 0.08865710560625815
$ws.Range("J19").Value = 0.3617992177314211
$ws.Range("K19").Value = 0.1140808344198175
$ws.Range("M19").Value = 0.01434159061277705
$ws.Range("N19").Value = 0.002607561929595828
$ws.Range("O19").Value = 0.06910039113428944
$ws.Range("S19").Value = 0.1544980443285528
